$d = $word.ActiveDocument

$d.Content.Find.Execute("2026-02-25 Wednesday", $true, $false, $false, $false, $false, $true, 1, $false, "2026-02-26 Thursday", 2) | Out-Null
$d.Content.Find.Execute("85+6=91", $true, $false, $false, $false, $false, $true, 1, $false, "69+2=71", 2) | Out-Null
$d.Content.Find.Execute("95-49=46", $true, $false, $false, $false, $false, $true, 1, $false, "88-1=87", 2) | Out-Null
$d.Content.Find.Execute("98-96=2", $true, $false, $false, $false, $false, $true, 1, $false, "68-37=31", 2) | Out-Null
$d.Content.Find.Execute("13+58=71", $true, $false, $false, $false, $false, $true, 1, $false, "99-13=86", 2) | Out-Null
$d.Content.Find.Execute("94-25=69", $true, $false, $false, $false, $false, $true, 1, $false, "3+60=63", 2) | Out-Null
$d.Content.Find.Execute("58+22=80", $true, $false, $false, $false, $false, $true, 1, $false, "63-0=63", 2) | Out-Null
$d.Content.Find.Execute("54+5=59", $true, $false, $false, $false, $false, $true, 1, $false, "88-27=61", 2) | Out-Null
$d.Content.Find.Execute("99-14=85", $true, $false, $false, $false, $false, $true, 1, $false, "11+76=87", 2) | Out-Null
$d.Content.Find.Execute("9+11=20", $true, $false, $false, $false, $false, $true, 1, $false, "90-66=24", 2) | Out-Null
$d.Content.Find.Execute("95-87=8", $true, $false, $false, $false, $false, $true, 1, $false, "86-21=65", 2) | Out-Null
$d.Content.Find.Execute("0+28=28", $true, $false, $false, $false, $false, $true, 1, $false, "96-65=31", 2) | Out-Null
$d.Content.Find.Execute("55+13=68", $true, $false, $false, $false, $false, $true, 1, $false, "61-19=42", 2) | Out-Null
$d.Content.Find.Execute("47+24=71", $true, $false, $false, $false, $false, $true, 1, $false, "77-24=53", 2) | Out-Null
$d.Content.Find.Execute("51+29=80", $true, $false, $false, $false, $false, $true, 1, $false, "22+22=44", 2) | Out-Null
$d.Content.Find.Execute("11-2=9", $true, $false, $false, $false, $false, $true, 1, $false, "4+30=34", 2) | Out-Null
$d.Content.Find.Execute("43+44=87", $true, $false, $false, $false, $false, $true, 1, $false, "89-7=82", 2) | Out-Null
$d.Content.Find.Execute("28+68=96", $true, $false, $false, $false, $false, $true, 1, $false, "55-2=53", 2) | Out-Null
$d.Content.Find.Execute("79-71=8", $true, $false, $false, $false, $false, $true, 1, $false, "87-12=75", 2) | Out-Null
$d.Content.Find.Execute("86-11=75", $true, $false, $false, $false, $false, $true, 1, $false, "9+66=75", 2) | Out-Null
$d.Content.Find.Execute("64-30=34", $true, $false, $false, $false, $false, $true, 1, $false, "57-44=13", 2) | Out-Null
$d.Content.Find.Execute("28+48=76", $true, $false, $false, $false, $false, $true, 1, $false, "35+58=93", 2) | Out-Null
$d.Content.Find.Execute("10+48=58", $true, $false, $false, $false, $false, $true, 1, $false, "82-51=31", 2) | Out-Null
$d.Content.Find.Execute("80-0=80", $true, $false, $false, $false, $false, $true, 1, $false, "58-48=10", 2) | Out-Null
$d.Content.Find.Execute("89-71=18", $true, $false, $false, $false, $false, $true, 1, $false, "83-39=44", 2) | Out-Null
$d.Content.Find.Execute("63+6=69", $true, $false, $false, $false, $false, $true, 1, $false, "88-8=80", 2) | Out-Null
$d.Content.Find.Execute("57-2=55", $true, $false, $false, $false, $false, $true, 1, $false, "55-8=47", 2) | Out-Null
$d.Content.Find.Execute("63-21=42", $true, $false, $false, $false, $false, $true, 1, $false, "30-22=8", 2) | Out-Null
$d.Content.Find.Execute("2+67=69", $true, $false, $false, $false, $false, $true, 1, $false, "20+2=22", 2) | Out-Null
$d.Content.Find.Execute("41+30=71", $true, $false, $false, $false, $false, $true, 1, $false, "14-7=7", 2) | Out-Null
$d.Content.Find.Execute("73+20=93", $true, $false, $false, $false, $false, $true, 1, $false, "95-13=82", 2) | Out-Null
$d.Content.Find.Execute("31+37=68", $true, $false, $false, $false, $false, $true, 1, $false, "59-35=24", 2) | Out-Null
$d.Content.Find.Execute("22+20=42", $true, $false, $false, $false, $false, $true, 1, $false, "83-39=44", 2) | Out-Null
$d.Content.Find.Execute("71-66=5", $true, $false, $false, $false, $false, $true, 1, $false, "1+12=13", 2) | Out-Null
$d.Content.Find.Execute("21+25=46", $true, $false, $false, $false, $false, $true, 1, $false, "54+29=83", 2) | Out-Null
$d.Content.Find.Execute("64-17=47", $true, $false, $false, $false, $false, $true, 1, $false, "20+43=63", 2) | Out-Null
$d.Content.Find.Execute("22+42=64", $true, $false, $false, $false, $false, $true, 1, $false, "83-13=70", 2) | Out-Null
$d.Content.Find.Execute("4+73=77", $true, $false, $false, $false, $false, $true, 1, $false, "83-81=2", 2) | Out-Null
$d.Content.Find.Execute("67-48=19", $true, $false, $false, $false, $false, $true, 1, $false, "20+39=59", 2) | Out-Null
$d.Content.Find.Execute("78-6=72", $true, $false, $false, $false, $false, $true, 1, $false, "1+68=69", 2) | Out-Null
$d.Content.Find.Execute("92-16=76", $true, $false, $false, $false, $false, $true, 1, $false, "14+30=44", 2) | Out-Null
$d.Content.Find.Execute("7+19=26", $true, $false, $false, $false, $false, $true, 1, $false, "95-33=62", 2) | Out-Null
$d.Content.Find.Execute("18+17=35", $true, $false, $false, $false, $false, $true, 1, $false, "54-0=54", 2) | Out-Null
$d.Content.Find.Execute("79-65=14", $true, $false, $false, $false, $false, $true, 1, $false, "27-10=17", 2) | Out-Null
$d.Content.Find.Execute("60-40=20", $true, $false, $false, $false, $false, $true, 1, $false, "4+14=18", 2) | Out-Null
$d.Content.Find.Execute("82-39=43", $true, $false, $false, $false, $false, $true, 1, $false, "44+50=94", 2) | Out-Null
$d.Content.Find.Execute("79-52=27", $true, $false, $false, $false, $false, $true, 1, $false, "19+45=64", 2) | Out-Null
$d.Content.Find.Execute("39+18=57", $true, $false, $false, $false, $false, $true, 1, $false, "14+20=34", 2) | Out-Null
$d.Content.Find.Execute("25-14=11", $true, $false, $false, $false, $false, $true, 1, $false, "42+26=68", 2) | Out-Null
$d.Content.Find.Execute("70-36=34", $true, $false, $false, $false, $false, $true, 1, $false, "14+9=23", 2) | Out-Null
$d.Content.Find.Execute("96-77=19", $true, $false, $false, $false, $false, $true, 1, $false, "65-15=50", 2) | Out-Null
$d.Content.Find.Execute("29+8=37", $true, $false, $false, $false, $false, $true, 1, $false, "99-31=68", 2) | Out-Null
$d.Content.Find.Execute("68-17=51", $true, $false, $false, $false, $false, $true, 1, $false, "79-78=1", 2) | Out-Null
$d.Content.Find.Execute("0+1=1", $true, $false, $false, $false, $false, $true, 1, $false, "87-3=84", 2) | Out-Null
$d.Content.Find.Execute("49-33=16", $true, $false, $false, $false, $false, $true, 1, $false, "78-15=63", 2) | Out-Null
$d.Content.Find.Execute("3+0=3", $true, $false, $false, $false, $false, $true, 1, $false, "29+10=39", 2) | Out-Null
$d.Content.Find.Execute("28-9=19", $true, $false, $false, $false, $false, $true, 1, $false, "96-66=30", 2) | Out-Null
$d.Content.Find.Execute("22+45=67", $true, $false, $false, $false, $false, $true, 1, $false, "6+62=68", 2) | Out-Null
$d.Content.Find.Execute("88-14=74", $true, $false, $false, $false, $false, $true, 1, $false, "19+75=94", 2) | Out-Null
$d.Content.Find.Execute("49+3=52", $true, $false, $false, $false, $false, $true, 1, $false, "7+26=33", 2) | Out-Null
$d.Content.Find.Execute("81-59=22", $true, $false, $false, $false, $false, $true, 1, $false, "60-32=28", 2) | Out-Null
$d.Content.Find.Execute("22+18=40", $true, $false, $false, $false, $false, $true, 1, $false, "75-8=67", 2) | Out-Null
$d.Content.Find.Execute("41+15=56", $true, $false, $false, $false, $false, $true, 1, $false, "19+44=63", 2) | Out-Null
$d.Content.Find.Execute("90+3=93", $true, $false, $false, $false, $false, $true, 1, $false, "59+15=74", 2) | Out-Null
$d.Content.Find.Execute("45+38=83", $true, $false, $false, $false, $false, $true, 1, $false, "91-48=43", 2) | Out-Null
$d.Content.Find.Execute("86-18=68", $true, $false, $false, $false, $false, $true, 1, $false, "32+34=66", 2) | Out-Null
$d.Content.Find.Execute("87-35=52", $true, $false, $false, $false, $false, $true, 1, $false, "4+0=4", 2) | Out-Null
$d.Content.Find.Execute("89-72=17", $true, $false, $false, $false, $false, $true, 1, $false, "31+68=99", 2) | Out-Null
$d.Content.Find.Execute("71-29=42", $true, $false, $false, $false, $false, $true, 1, $false, "98-40=58", 2) | Out-Null
$d.Content.Find.Execute("58-1=57", $true, $false, $false, $false, $false, $true, 1, $false, "73-11=62", 2) | Out-Null
$d.Content.Find.Execute("3+75=78", $true, $false, $false, $false, $false, $true, 1, $false, "30-23=7", 2) | Out-Null
$d.Content.Find.Execute("45+9=54", $true, $false, $false, $false, $false, $true, 1, $false, "57-41=16", 2) | Out-Null
$d.Content.Find.Execute("7+14=21", $true, $false, $false, $false, $false, $true, 1, $false, "72+19=91", 2) | Out-Null
$d.Content.Find.Execute("42+14=56", $true, $false, $false, $false, $false, $true, 1, $false, "29-25=4", 2) | Out-Null
$d.Content.Find.Execute("14+49=63", $true, $false, $false, $false, $false, $true, 1, $false, "93+3=96", 2) | Out-Null
$d.Content.Find.Execute("69-5=64", $true, $false, $false, $false, $false, $true, 1, $false, "8+53=61", 2) | Out-Null
$d.Content.Find.Execute("74-44=30", $true, $false, $false, $false, $false, $true, 1, $false, "86-61=25", 2) | Out-Null
$d.Content.Find.Execute("35-8=27", $true, $false, $false, $false, $false, $true, 1, $false, "2+9=11", 2) | Out-Null
$d.Content.Find.Execute("84-80=4", $true, $false, $false, $false, $false, $true, 1, $false, "93-50=43", 2) | Out-Null
$d.Content.Find.Execute("57-36=21", $true, $false, $false, $false, $false, $true, 1, $false, "86-28=58", 2) | Out-Null
$d.Content.Find.Execute("12+34=46", $true, $false, $false, $false, $false, $true, 1, $false, "12+7=19", 2) | Out-Null
$d.Content.Find.Execute("5+60=65", $true, $false, $false, $false, $false, $true, 1, $false, "16+46=62", 2) | Out-Null
$d.Content.Find.Execute("19-11=8", $true, $false, $false, $false, $false, $true, 1, $false, "76-50=26", 2) | Out-Null
$d.Content.Find.Execute("16+55=71", $true, $false, $false, $false, $false, $true, 1, $false, "81-41=40", 2) | Out-Null
$d.Content.Find.Execute("45+30=75", $true, $false, $false, $false, $false, $true, 1, $false, "3+78=81", 2) | Out-Null
$d.Content.Find.Execute("37+39=76", $true, $false, $false, $false, $false, $true, 1, $false, "64-41=23", 2) | Out-Null
$d.Content.Find.Execute("51-43=8", $true, $false, $false, $false, $false, $true, 1, $false, "4+21=25", 2) | Out-Null
$d.Content.Find.Execute("11+4=15", $true, $false, $false, $false, $false, $true, 1, $false, "84-17=67", 2) | Out-Null
$d.Content.Find.Execute("22-9=13", $true, $false, $false, $false, $false, $true, 1, $false, "38+34=72", 2) | Out-Null
$d.Content.Find.Execute("88-54=34", $true, $false, $false, $false, $false, $true, 1, $false, "9+66=75", 2) | Out-Null
$d.Content.Find.Execute("75-45=30", $true, $false, $false, $false, $false, $true, 1, $false, "72-16=56", 2) | Out-Null
$d.Content.Find.Execute("19+29=48", $true, $false, $false, $false, $false, $true, 1, $false, "15+50=65", 2) | Out-Null
$d.Content.Find.Execute("86-81=5", $true, $false, $false, $false, $false, $true, 1, $false, "20+40=60", 2) | Out-Null
$d.Content.Find.Execute("41-5=36", $true, $false, $false, $false, $false, $true, 1, $false, "65+10=75", 2) | Out-Null
$d.Content.Find.Execute("38+60=98", $true, $false, $false, $false, $false, $true, 1, $false, "94-58=36", 2) | Out-Null
$d.Content.Find.Execute("29-28=1", $true, $false, $false, $false, $false, $true, 1, $false, "54-35=19", 2) | Out-Null
$d.Content.Find.Execute("46-31=15", $true, $false, $false, $false, $false, $true, 1, $false, "45-21=24", 2) | Out-Null
$d.Content.Find.Execute("35+52=87", $true, $false, $false, $false, $false, $true, 1, $false, "63-46=17", 2) | Out-Null
$d.Content.Find.Execute("97-95=2", $true, $false, $false, $false, $false, $true, 1, $false, "71-42=29", 2) | Out-Null
$d.Content.Find.Execute("73-15=58", $true, $false, $false, $false, $false, $true, 1, $false, "40+34=74", 2) | Out-Null
$d.Content.Find.Execute("17+66=83", $true, $false, $false, $false, $false, $true, 1, $false, "43-27=16", 2) | Out-Null

Write-Output "Replacements complete"